{"js": "const replacements = [\n  { find: \"2025-11-17 Monday\", replace: \"2025-11-18 Tuesday\" },\n  { find: \"45\u00f79=\", replace: \"85\u00f74=\" },\n  { find: \"60\u00f74=\", replace: \"33\u00f77=\" },\n  { find: \"45\u00f75=\", replace: \"75\u00f76=\" },\n  { find: \"58\u00f74=\", replace: \"90\u00f79=\" },\n  { find: \"48\u00f72=\", replace: \"87\u00f78=\" },\n  { find: \"26\u00f79=\", replace: \"47\u00f78=\" },\n  { find: \"73\u00f75=\", replace: \"43\u00f73=\" },\n  { find: \"84\u00f73=\", replace: \"97\u00f76=\" },\n  { find: \"30\u00f79=\", replace: \"77\u00f74=\" },\n  { find: \"43\u00f74=\", replace: \"77\u00f75=\" },\n  { find: \"67\u00f74=\", replace: \"50\u00f72=\" },\n  { find: \"33\u00f72=\", replace: \"97\u00f79=\" },\n  { find: \"78\u00f76=\", replace: \"89\u00f78=\" },\n  { find: \"64\u00f76=\", replace: \"91\u00f77=\" },\n  { find: \"93\u00f76=\", replace: \"48\u00f74=\" },\n  { find: \"88\u00f79=\", replace: \"17\u00f76=\" },\n  { find: \"68\u00f76=\", replace: \"12\u00f76=\" },\n  { find: \"58\u00f73=\", replace: \"59\u00f74=\" },\n  { find: \"76\u00f78=\", replace: \"36\u00f74=\" },\n  { find: \"35\u00f79=\", replace: \"43\u00f73=\" },\n  { find: \"10\u00f79=\", replace: \"26\u00f77=\" },\n  { find: \"86\u00f75=\", replace: \"97\u00f77=\" },\n  { find: \"83\u00f78=\", replace: \"32\u00f78=\" },\n  { find: \"18\u00f72=\", replace: \"77\u00f79=\" },\n  { find: \"90\u00f75=\", replace: \"82\u00f78=\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$wdReplaceAll = 2\n\n$replacements = @(\n    @{ Find = '2025-11-17 Monday'; Replace = '2025-11-18 Tuesday' }\n    @{ Find = '45\u00f79='; Replace = '85\u00f74=' }\n    @{ Find = '60\u00f74='; Replace = '33\u00f77=' }\n    @{ Find = '45\u00f75='; Replace = '75\u00f76=' }\n    @{ Find = '58\u00f74='; Replace = '90\u00f79=' }\n    @{ Find = '48\u00f72='; Replace = '87\u00f78=' }\n    @{ Find = '26\u00f79='; Replace = '47\u00f78=' }\n    @{ Find = '73\u00f75='; Replace = '43\u00f73=' }\n    @{ Find = '84\u00f73='; Replace = '97\u00f76=' }\n    @{ Find = '30\u00f79='; Replace = '77\u00f74=' }\n    @{ Find = '43\u00f74='; Replace = '77\u00f75=' }\n    @{ Find = '67\u00f74='; Replace = '50\u00f72=' }\n    @{ Find = '33\u00f72='; Replace = '97\u00f79=' }\n    @{ Find = '78\u00f76='; Replace = '89\u00f78=' }\n    @{ Find = '64\u00f76='; Replace = '91\u00f77=' }\n    @{ Find = '93\u00f76='; Replace = '48\u00f74=' }\n    @{ Find = '88\u00f79='; Replace = '17\u00f76=' }\n    @{ Find = '68\u00f76='; Replace = '12\u00f76=' }\n    @{ Find = '58\u00f73='; Replace = '59\u00f74=' }\n    @{ Find = '76\u00f78='; Replace = '36\u00f74=' }\n    @{ Find = '35\u00f79='; Replace = '43\u00f73=' }\n    @{ Find = '10\u00f79='; Replace = '26\u00f77=' }\n    @{ Find = '86\u00f75='; Replace = '97\u00f77=' }\n    @{ Find = '83\u00f78='; Replace = '32\u00f78=' }\n    @{ Find = '18\u00f72='; Replace = '77\u00f79=' }\n    @{ Find = '90\u00f75='; Replace = '82\u00f78=' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, $wdReplaceAll)\n}\n"}
